# Updates cryptos list: Price (D) and Volume(1h) (E) columns for rows 2-51
# matching the commit "Updated cryptos list on Sun Sep 10 17:00:21 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a pure number by Excel
# (losing formatting like trailing zeros, e.g. "61.20" -> 61.2). Force them to
# plain text first, write the value, then restore the default "Normal" style so
# no stray number-format/style is left behind on the cell.
$textForceCells = @("D5","D6","D8","D9","D10","D11","D15","D18","D20","D22","D23","D25","D27","D28","D29","D31","D33","D35","D37","D39","D41","D42","D44","D46","D47","D48")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.031.12"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.625.56"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "214.01"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").Value = "0.0618"
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("D10").Value = "18.18"
$ws.Range("E10").Value = "  -7.42%  "
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").Value = "1.852.76"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.643.22"
$ws.Range("E13").Value = "  -3.85%  "
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").Value = "0.522"
$ws.Range("E15").Value = "  -4.07%  "
$ws.Range("D16").Value = "25.995.52"
$ws.Range("D17").Value = "0.0₃0739"
$ws.Range("E17").Value = "  -3.26%  "
$ws.Range("D18").Value = "61.20"
$ws.Range("E18").Value = "  -3.51%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "189.57"
$ws.Range("E20").Value = "  -3.15%  "
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").Value = "9.56"
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("D23").Value = "6.06"
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "143.82"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "1.77"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").Value = "6.71"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("D29").Value = "15.13"
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").Value = "0.0481"
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("D33").Value = "3.11"
$ws.Range("E33").Value = "  -5.73%  "
$ws.Range("E34").Value = "  -2.19%  "
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  -2.92%  "
$ws.Range("D36").Value = "1.127.59"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").Value = "0.849"
$ws.Range("E37").Value = "  -6.25%  "
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").Value = "0.515"
$ws.Range("E39").Value = "  -4.88%  "
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("D41").Value = "97.97"
$ws.Range("E41").Value = "  -1.38%  "
$ws.Range("D42").Value = "0.774"
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("D43").Value = "1.765.12"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "5.17"
$ws.Range("E44").Value = "  -5.76%  "
$ws.Range("D45").Value = "0.0₆0109"
$ws.Range("E45").Value = "  -7.39%  "
$ws.Range("D46").Value = "54.48"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").Value = "0.0526"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "1.48"
$ws.Range("E48").Value = "  +0.39%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("E51").Value = "  -3.99%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
